$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Date in A1 bumped by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Updated prices in column D for rows 31-38
$ws.Range("D31").Value = 4359.3
$ws.Range("D32").Value = 3888.5
$ws.Range("D33").Value = 3708.1
$ws.Range("D34").Value = 6315.1
$ws.Range("D35").Value = 4712.4
$ws.Range("D36").Value = 4441.8
$ws.Range("D37").Value = 4162.4
$ws.Range("D38").Value = 6581.3
